# Apply cryptos.xlsx price/volume updates (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.728.77"
$ws.Range("E2").Value = "  +1.57%  "

$ws.Range("D3").Value = "1.646.90"
$ws.Range("E3").Value = "  +0.58%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").Value = "213.27"
$ws.Range("E5").Value = "  +0.65%  "

$ws.Range("D6").Value = "0.528"
$ws.Range("E6").Value = "  -0.25%  "

$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("D8").Value = "23.43"
$ws.Range("E8").Value = "  +1.30%  "

$ws.Range("E9").Value = "  +1.08%  "

$ws.Range("D10").Value = "0.0614"
$ws.Range("E10").Value = "  +1.05%  "

$ws.Range("D11").Value = "0.0888"
$ws.Range("E11").Value = "  -0.05%  "

$ws.Range("D12").Value = "1.881.40"
$ws.Range("E12").Value = "  +0.60%  "

$ws.Range("D13").Value = "1.644.72"
$ws.Range("E13").Value = "  +0.21%  "

$ws.Range("E14").Value = "  +0.98%  "

$ws.Range("D15").Value = "'0.560"
$ws.Range("E15").Value = "  +0.69%  "

$ws.Range("D16").Value = "64.54"
$ws.Range("E16").Value = "  +0.67%  "

$ws.Range("D17").Value = "27.731.20"
$ws.Range("E17").Value = "  +1.64%  "

$ws.Range("D18").Value = "232.37"
$ws.Range("E18").Value = "  +2.42%  "

$ws.Range("D19").Value = "0.0₃0725"
$ws.Range("E19").Value = "  +1.16%  "

$ws.Range("D20").Value = "7.61"
$ws.Range("E20").Value = "  +3.82%  "

$ws.Range("E21").Value = "  +0.09%  "

$ws.Range("D22").Value = "4.29"
$ws.Range("E22").Value = "  +0.36%  "

$ws.Range("D23").Value = "10.13"
$ws.Range("E23").Value = "  +9.92%  "

$ws.Range("D24").Value = "1.95"
$ws.Range("E24").Value = "  -3.27%  "

$ws.Range("D25").Value = "149.54"
$ws.Range("E25").Value = "  +1.77%  "

$ws.Range("D26").Value = "6.94"
$ws.Range("E26").Value = "  +0.53%  "

$ws.Range("D27").Value = "0.111"
$ws.Range("E27").Value = "  -1.10%  "

$ws.Range("D28").Value = "15.68"
$ws.Range("E28").Value = "  +1.48%  "

$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.09%  "

$ws.Range("E30").Value = "  +0.85%  "

$ws.Range("D31").Value = "0.0487"
$ws.Range("E31").Value = "  +0.52%  "

$ws.Range("D32").Value = "3.29"
$ws.Range("E32").Value = "  +1.17%  "

$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "3.15"
$ws.Range("E33").Value = "  +2.54%  "

$ws.Range("B34").Value = "Maker"
$ws.Range("C34").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D34").Value = "1.433.45"
$ws.Range("E34").Value = "  +3.09%  "

$ws.Range("D35").Value = "1.59"
$ws.Range("E35").Value = "  +2.43%  "

$ws.Range("E36").Value = "  -1.26%  "

$ws.Range("D37").Value = "0.568"
$ws.Range("E37").Value = "  +2.81%  "

$ws.Range("D38").Value = "0.883"
$ws.Range("E38").Value = "  +1.23%  "

$ws.Range("D39").Value = "0.0167"
$ws.Range("E39").Value = "  +1.58%  "

$ws.Range("D40").Value = "0.881"
$ws.Range("E40").Value = "  +12.07%  "

$ws.Range("D41").Value = "1.04"
$ws.Range("E41").Value = "  +1.19%  "

$ws.Range("E42").Value = "  +0.12%  "

$ws.Range("D43").Value = "5.58"
$ws.Range("E43").Value = "  +2.92%  "

$ws.Range("D44").Value = "66.82"
$ws.Range("E44").Value = "  +4.90%  "

$ws.Range("B45").Value = "mCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$ws.Range("D45").Value = "2.46"
$ws.Range("E45").Value = "  -0.87%  "

$ws.Range("B46").Value = "MXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D46").Value = "2.26"
$ws.Range("E46").Value = "  +2.12%  "

$ws.Range("D47").Value = "1.790.66"
$ws.Range("E47").Value = "  +0.51%  "

$ws.Range("E48").Value = "  +6.81%  "

$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₆0108"
$ws.Range("E49").Value = "  +2.89%  "

$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "85.24"
$ws.Range("E50").Value = "  -1.81%  "

$ws.Range("D51").Value = "'0.0990"
$ws.Range("E51").Value = "  +1.30%  "
